# Adding (better) commentary to code ++ Prepwork for time-shift capabilities. ++ Bugfixes
#
# Prepwork for time-shift capabilities: add three new ColumnHeader entries
# (clh_ModifyDate / clh_CreateDate / clh_OffsetTime) with their display-name
# counterparts (Modify Date / Create Date / Time Offset) to the
# "ColumnHeader" lookup sheet, then leave that sheet active/selected as it
# was left in the authoring session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnHeader")

# New internal (code) names go in column A, first, so the shared-string
# table picks them up before the human-readable labels.
$ws.Cells.Item(33, 1).Value2 = "clh_ModifyDate"
$ws.Cells.Item(34, 1).Value2 = "clh_CreateDate"
$ws.Cells.Item(35, 1).Value2 = "clh_OffsetTime"

# Human readable labels go in column C.
$ws.Cells.Item(33, 3).Value2 = "Modify Date"
$ws.Cells.Item(34, 3).Value2 = "Create Date"
$ws.Cells.Item(35, 3).Value2 = "Time Offset"

# Leave the ColumnHeader tab as the active/selected one, with the cell just
# below the newly-added rows selected.
$ws.Activate()
$ws.Range("C36").Select()
